$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two rows (Uzbek, Vietnamese) - rows 22 and 23
$ws.Range("A22:B23").Delete() | Out-Null

# Re-order the remaining 20 language rows (rows 2-21) by descending value
$data = @(
    @("English", 27.10231750593171),
    @("Spanish", 8.210106001566238),
    @("Chinese", 7.990769160021653),
    @("Japanese", 7.284845128816129),
    @("German", 6.00564843394439),
    @("Arabic", 4.854315317304337),
    @("Portuguese", 3.788016098813561),
    @("French", 3.556750962104147),
    @("Italian", 3.383359981502351),
    @("Russian", 2.934035637401106),
    @("Malay-Indonesian", 2.567766416602798),
    @("Dutch", 1.647057393752311),
    @("Turkish", 1.453072909350366),
    @("Persian", 1.41649117189301),
    @("Korean", 1.393184364900861),
    @("Thai", 0.9023102007134319),
    @("Polish", 0.8882295149114708),
    @("Urdu", 0.8051632623911761),
    @("Swedish", 0.5130925337219181),
    @("Bengali", 0.4224161223718882)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
